# Updated cryptos list - applies the per-cell text updates described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($Cell, [string]$NewValue)
    # Force the cell to be treated as text so that values such as
    # "536.17" or "62.901.97" are stored verbatim instead of being
    # auto-converted into (imprecise) floating point numbers.
    $Cell.NumberFormat = "@"
    $Cell.Value = $NewValue
    $Cell.ClearFormats()
}

Set-TextCell $ws.Range("D2") "62.901.97"
Set-TextCell $ws.Range("E2") "  -1.11%  "
Set-TextCell $ws.Range("D3") "3.069.42"
Set-TextCell $ws.Range("E3") "  -0.45%  "
Set-TextCell $ws.Range("E4") "  -0.09%  "
Set-TextCell $ws.Range("D5") "536.17"
Set-TextCell $ws.Range("E5") "  -3.22%  "
Set-TextCell $ws.Range("D6") "133.44"
Set-TextCell $ws.Range("E6") "  -2.12%  "
Set-TextCell $ws.Range("D7") "0.999"
Set-TextCell $ws.Range("E7") "  -0.07%  "
Set-TextCell $ws.Range("D8") "3.066.41"
Set-TextCell $ws.Range("E8") "  -0.23%  "
Set-TextCell $ws.Range("D9") "0.494"
Set-TextCell $ws.Range("E9") "  +0.72%  "
Set-TextCell $ws.Range("E10") "  -1.42%  "
Set-TextCell $ws.Range("D11") "6.22"
Set-TextCell $ws.Range("E11") "  -5.61%  "
Set-TextCell $ws.Range("D12") "0.453"
Set-TextCell $ws.Range("E12") "  +1.36%  "
Set-TextCell $ws.Range("E13") "  +3.99%  "
Set-TextCell $ws.Range("D14") "34.16"
Set-TextCell $ws.Range("E14") "  -1.76%  "
Set-TextCell $ws.Range("D15") "3.561.10"
Set-TextCell $ws.Range("E15") "  -0.57%  "
Set-TextCell $ws.Range("D16") "62.843.12"
Set-TextCell $ws.Range("E16") "  -1.40%  "
Set-TextCell $ws.Range("E17") "  -0.13%  "
Set-TextCell $ws.Range("D18") "3.065.37"
Set-TextCell $ws.Range("E18") "  -0.63%  "
Set-TextCell $ws.Range("D19") "6.63"
Set-TextCell $ws.Range("E19") "  +1.14%  "
Set-TextCell $ws.Range("D20") "481.37"
Set-TextCell $ws.Range("E20") "  -4.17%  "
Set-TextCell $ws.Range("D21") "13.31"
Set-TextCell $ws.Range("E21") "  -1.07%  "
Set-TextCell $ws.Range("D22") "0.693"
Set-TextCell $ws.Range("E22") "  -0.44%  "
Set-TextCell $ws.Range("D23") "7.12"
Set-TextCell $ws.Range("E23") "  -0.38%  "
Set-TextCell $ws.Range("D24") "79.01"
Set-TextCell $ws.Range("E24") "  +2.92%  "
Set-TextCell $ws.Range("E25") "  -1.58%  "
Set-TextCell $ws.Range("D26") "0.999"
Set-TextCell $ws.Range("E26") "  +0.16%  "
Set-TextCell $ws.Range("E27") "  -2.34%  "
Set-TextCell $ws.Range("D28") "8.07"
Set-TextCell $ws.Range("E28") "  -1.48%  "
Set-TextCell $ws.Range("D29") "0.999"
Set-TextCell $ws.Range("E29") "  -0.09%  "
Set-TextCell $ws.Range("D30") "25.98"
Set-TextCell $ws.Range("E30") "  +0.12%  "
Set-TextCell $ws.Range("D31") "1.86"
Set-TextCell $ws.Range("E31") "  -7.56%  "
Set-TextCell $ws.Range("E32") "  +1.04%  "
Set-TextCell $ws.Range("B33") "OKB"
Set-TextCell $ws.Range("C33") "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextCell $ws.Range("D33") "56.84"
Set-TextCell $ws.Range("E33") "  -1.52%  "
Set-TextCell $ws.Range("B34") "Stacks"
Set-TextCell $ws.Range("C34") "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextCell $ws.Range("D34") "2.35"
Set-TextCell $ws.Range("E34") "  -6.16%  "
Set-TextCell $ws.Range("D35") "5.35"
Set-TextCell $ws.Range("E35") "  +5.06%  "
Set-TextCell $ws.Range("D36") "5.97"
Set-TextCell $ws.Range("E36") "  +3.02%  "
Set-TextCell $ws.Range("D37") "477.71"
Set-TextCell $ws.Range("E37") "  -9.62%  "
Set-TextCell $ws.Range("B38") "Maker"
Set-TextCell $ws.Range("C38") "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextCell $ws.Range("D38") "3.084.54"
Set-TextCell $ws.Range("E38") "  +1.33%  "
Set-TextCell $ws.Range("B39") "VeChain"
Set-TextCell $ws.Range("C39") "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextCell $ws.Range("D39") "0.0392"
Set-TextCell $ws.Range("E39") "  -3.99%  "
Set-TextCell $ws.Range("D40") "0.0795"
Set-TextCell $ws.Range("E40") "  +1.09%  "
Set-TextCell $ws.Range("E41") "  -1.40%  "
Set-TextCell $ws.Range("E42") "  +1.11%  "
Set-TextCell $ws.Range("E43") "  +4.12%  "
Set-TextCell $ws.Range("D44") "0.252"
Set-TextCell $ws.Range("E44") "  +1.36%  "
Set-TextCell $ws.Range("E46") "  +10.43%  "
Set-TextCell $ws.Range("D47") "121.29"
Set-TextCell $ws.Range("E47") "  -0.56%  "
Set-TextCell $ws.Range("E48") "  -1.76%  "
Set-TextCell $ws.Range("D49") "24.41"
Set-TextCell $ws.Range("E49") "  +2.27%  "
Set-TextCell $ws.Range("E50") "  +2.28%  "
Set-TextCell $ws.Range("D51") "2.29"
Set-TextCell $ws.Range("E51") "  +5.32%  "
